$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B2").Value = 45.48165850083172
$ws.Range("B3").Value = 29.14359298898636
$ws.Range("B4").Value = 10.2492150563531
$ws.Range("B5").Value = 6.737806404501392
$ws.Range("B6").Value = 4.128183617122676
$ws.Range("B7").Value = 2.266875590221579
$ws.Range("B8").Value = 1.992667841983179
